$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# --- ALC ---
$ws_ALC.Range("H33").Value = 379.8857
$ws_ALC.Range("J33").Value = 513.2857
$ws_ALC.Range("L33").Value = 513.2857
$ws_ALC.Range("N33").Value = -971.2857
$ws_ALC.Range("H64").Value = 4362.5
$ws_ALC.Range("I64").Value = 4745
$ws_ALC.Range("J64").Value = 3980
$ws_ALC.Range("K64").Value = 4745
$ws_ALC.Range("L64").Value = 3980
$ws_ALC.Range("M64").Value = -4497
$ws_ALC.Range("N64").Value = -4476
$ws_ALC.Range("H67").Value = 4362.5
$ws_ALC.Range("I67").Value = 4745
$ws_ALC.Range("J67").Value = 3980
$ws_ALC.Range("K67").Value = 4745
$ws_ALC.Range("L67").Value = 3980
$ws_ALC.Range("M67").Value = -3887
$ws_ALC.Range("N67").Value = -5696
$ws_ALC.Range("H121").Value = 1498.3334
$ws_ALC.Range("I121").Value = 1500
$ws_ALC.Range("J121").Value = 1497.5
$ws_ALC.Range("K121").Value = 4500
$ws_ALC.Range("L121").Value = 4492.5
$ws_ALC.Range("M121").Value = -2753
$ws_ALC.Range("N121").Value = -7986.5
$ws_ALC.Range("H137").Value = 2062.739
$ws_ALC.Range("I137").Value = 1463.8
$ws_ALC.Range("J137").Value = 3185.75
$ws_ALC.Range("K137").Value = 4391.4
$ws_ALC.Range("L137").Value = 9557.25
$ws_ALC.Range("M137").Value = -1841.4
$ws_ALC.Range("N137").Value = -14657.25

# --- ARM ---
$ws_ARM.Range("H61").Value = 1181.5
$ws_ARM.Range("I61").Value = 904
$ws_ARM.Range("K61").Value = 904
$ws_ARM.Range("M61").Value = -692
$ws_ARM.Range("H74").Value = 2027.4706
$ws_ARM.Range("I74").Value = 1737.6428
$ws_ARM.Range("J74").Value = 3380
$ws_ARM.Range("K74").Value = 1737.6428
$ws_ARM.Range("L74").Value = 3380
$ws_ARM.Range("M74").Value = -863.6428000000001
$ws_ARM.Range("N74").Value = -5128
$ws_ARM.Range("H77").Value = 2027.4706
$ws_ARM.Range("I77").Value = 1737.6428
$ws_ARM.Range("J77").Value = 3380
$ws_ARM.Range("K77").Value = 8688.214
$ws_ARM.Range("L77").Value = 16900
$ws_ARM.Range("M77").Value = -4320.214
$ws_ARM.Range("N77").Value = -25636
$ws_ARM.Range("H102").Value = 16669137
$ws_ARM.Range("I102").Value = 20835670
$ws_ARM.Range("K102").Value = 20835670
$ws_ARM.Range("M102").Value = -20834048
$ws_ARM.Range("H132").Value = 2327.6
$ws_ARM.Range("I132").Value = 1855.2941
$ws_ARM.Range("K132").Value = 5565.8823
$ws_ARM.Range("M132").Value = -3035.8823
$ws_ARM.Range("H136").Value = 1181.5
$ws_ARM.Range("I136").Value = 904
$ws_ARM.Range("K136").Value = 2712
$ws_ARM.Range("M136").Value = -162

# --- BSM ---
$ws_BSM.Range("H110").Value = 49999.75
$ws_BSM.Range("J110").Value = 49999.75
$ws_BSM.Range("L110").Value = 49999.75
$ws_BSM.Range("N110").Value = -58179.75
$ws_BSM.Range("H134").Value = 5140.5
$ws_BSM.Range("I134").Value = 985.5417
$ws_BSM.Range("K134").Value = 2956.6251
$ws_BSM.Range("M134").Value = -421.6251000000002

# --- CRP ---
$ws_CRP.Range("H16").Value = 71430050
$ws_CRP.Range("I16").Value = 111112520
$ws_CRP.Range("J16").Value = 1613.4
$ws_CRP.Range("K16").Value = 111112520
$ws_CRP.Range("L16").Value = 1613.4
$ws_CRP.Range("M16").Value = -111112233
$ws_CRP.Range("N16").Value = -2187.4
$ws_CRP.Range("H18").Value = 47400
$ws_CRP.Range("J18").Value = 47400
$ws_CRP.Range("L18").Value = 47400
$ws_CRP.Range("N18").Value = -47860
$ws_CRP.Range("H31").Value = 1299.3846
$ws_CRP.Range("I31").Value = 1260
$ws_CRP.Range("J31").Value = 1362.4
$ws_CRP.Range("K31").Value = 1260
$ws_CRP.Range("L31").Value = 1362.4
$ws_CRP.Range("M31").Value = -965
$ws_CRP.Range("N31").Value = -1952.4
$ws_CRP.Range("H34").Value = 1299.3846
$ws_CRP.Range("I34").Value = 1260
$ws_CRP.Range("J34").Value = 1362.4
$ws_CRP.Range("K34").Value = 1260
$ws_CRP.Range("L34").Value = 1362.4
$ws_CRP.Range("M34").Value = -1058
$ws_CRP.Range("N34").Value = -1766.4
$ws_CRP.Range("H86").Value = 3051263.5
$ws_CRP.Range("I86").Value = 5566031
$ws_CRP.Range("J86").Value = 33542.7
$ws_CRP.Range("K86").Value = 5566031
$ws_CRP.Range("L86").Value = 33542.7
$ws_CRP.Range("M86").Value = -5564908
$ws_CRP.Range("N86").Value = -35788.7
$ws_CRP.Range("H89").Value = 3051263.5
$ws_CRP.Range("I89").Value = 5566031
$ws_CRP.Range("J89").Value = 33542.7
$ws_CRP.Range("K89").Value = 27830155
$ws_CRP.Range("L89").Value = 167713.5
$ws_CRP.Range("M89").Value = -27824539
$ws_CRP.Range("N89").Value = -178945.5
$ws_CRP.Range("H105").Value = 727.5
$ws_CRP.Range("I105").Value = 655
$ws_CRP.Range("J105").Value = 800
$ws_CRP.Range("K105").Value = 655
$ws_CRP.Range("L105").Value = 800
$ws_CRP.Range("M105").Value = 1092
$ws_CRP.Range("N105").Value = -4294
$ws_CRP.Range("H107").Value = 740.95
$ws_CRP.Range("I107").Value = 369.9375
$ws_CRP.Range("J107").Value = 2225
$ws_CRP.Range("K107").Value = 369.9375
$ws_CRP.Range("L107").Value = 2225
$ws_CRP.Range("M107").Value = 1550.0625
$ws_CRP.Range("N107").Value = -6065
$ws_CRP.Range("H113").Value = 71430050
$ws_CRP.Range("I113").Value = 111112520
$ws_CRP.Range("J113").Value = 1613.4
$ws_CRP.Range("K113").Value = 111112520
$ws_CRP.Range("L113").Value = 1613.4
$ws_CRP.Range("M113").Value = -111110350
$ws_CRP.Range("N113").Value = -5953.4
$ws_CRP.Range("H132").Value = 2121.8262
$ws_CRP.Range("I132").Value = 1801.091
$ws_CRP.Range("K132").Value = 5403.272999999999
$ws_CRP.Range("M132").Value = -2873.272999999999
$ws_CRP.Range("H134").Value = 2072.875
$ws_CRP.Range("I134").Value = 2033.409
$ws_CRP.Range("J134").Value = 2507
$ws_CRP.Range("K134").Value = 6100.227000000001
$ws_CRP.Range("L134").Value = 7521
$ws_CRP.Range("M134").Value = -3565.227000000001
$ws_CRP.Range("N134").Value = -12591
$ws_CRP.Range("H141").Value = 874192.9
$ws_CRP.Range("J141").Value = 874192.9
$ws_CRP.Range("L141").Value = 874192.9
$ws_CRP.Range("N141").Value = -884552.9

# --- CUL ---
$ws_CUL.Range("H20").Value = 200
$ws_CUL.Range("J20").Value = 0
$ws_CUL.Range("L20").Value = 0
$ws_CUL.Range("N20").ClearContents()
$ws_CUL.Range("H22").Value = 3000
$ws_CUL.Range("J22").Value = 5000
$ws_CUL.Range("L22").Value = 15000
$ws_CUL.Range("N22").Value = -15338
$ws_CUL.Range("H27").Value = 3000
$ws_CUL.Range("J27").Value = 5000
$ws_CUL.Range("L27").Value = 15000
$ws_CUL.Range("N27").Value = -15204
$ws_CUL.Range("H39").Value = 3228.4285
$ws_CUL.Range("J39").Value = 3441.5
$ws_CUL.Range("L39").Value = 10324.5
$ws_CUL.Range("N39").Value = -10912.5

# --- GSM ---
$ws_GSM.Range("H10").Value = 5000334.5
$ws_GSM.Range("I10").Value = 5000334.5
$ws_GSM.Range("K10").Value = 5000334.5
$ws_GSM.Range("M10").Value = -5000165.5
$ws_GSM.Range("H113").Value = 1219.5
$ws_GSM.Range("I113").Value = 1412
$ws_GSM.Range("J113").Value = 1112.5555
$ws_GSM.Range("K113").Value = 1412
$ws_GSM.Range("L113").Value = 1112.5555
$ws_GSM.Range("M113").Value = 758
$ws_GSM.Range("N113").Value = -5452.5555
$ws_GSM.Range("H126").Value = 2151.0833
$ws_GSM.Range("I126").Value = 1824.875
$ws_GSM.Range("J126").Value = 2803.5
$ws_GSM.Range("K126").Value = 5474.625
$ws_GSM.Range("L126").Value = 8410.5
$ws_GSM.Range("M126").Value = -3004.625
$ws_GSM.Range("N126").Value = -13350.5
$ws_GSM.Range("H132").Value = 4630.143
$ws_GSM.Range("I132").Value = 5238.75
$ws_GSM.Range("K132").Value = 15716.25
$ws_GSM.Range("M132").Value = -13186.25

# --- LTW ---
$ws_LTW.Range("H7").Value = 1877
$ws_LTW.Range("I7").Value = 1596.25
$ws_LTW.Range("J7").Value = 3000
$ws_LTW.Range("K7").Value = 1596.25
$ws_LTW.Range("L7").Value = 3000
$ws_LTW.Range("M7").Value = -1484.25
$ws_LTW.Range("N7").Value = -3224
$ws_LTW.Range("H61").Value = 1136.8182
$ws_LTW.Range("I61").Value = 1055.5555
$ws_LTW.Range("J61").Value = 1502.5
$ws_LTW.Range("K61").Value = 1055.5555
$ws_LTW.Range("L61").Value = 1502.5
$ws_LTW.Range("M61").Value = -853.5554999999999
$ws_LTW.Range("N61").Value = -1906.5
$ws_LTW.Range("H113").Value = 1136.8182
$ws_LTW.Range("I113").Value = 1055.5555
$ws_LTW.Range("J113").Value = 1502.5
$ws_LTW.Range("K113").Value = 1055.5555
$ws_LTW.Range("L113").Value = 1502.5
$ws_LTW.Range("M113").Value = 1114.4445
$ws_LTW.Range("N113").Value = -5842.5
$ws_LTW.Range("H122").Value = 50009140
$ws_LTW.Range("I122").Value = 83334630
$ws_LTW.Range("J122").Value = 20900
$ws_LTW.Range("K122").Value = 250003890
$ws_LTW.Range("L122").Value = 62700
$ws_LTW.Range("M122").Value = -250001440
$ws_LTW.Range("N122").Value = -67600
$ws_LTW.Range("H126").Value = 1877
$ws_LTW.Range("I126").Value = 1596.25
$ws_LTW.Range("J126").Value = 3000
$ws_LTW.Range("K126").Value = 4788.75
$ws_LTW.Range("L126").Value = 9000
$ws_LTW.Range("M126").Value = -2318.75
$ws_LTW.Range("N126").Value = -13940
$ws_LTW.Range("H132").Value = 76793
$ws_LTW.Range("I132").Value = 24660.4
$ws_LTW.Range("J132").Value = 102859.3
$ws_LTW.Range("K132").Value = 73981.20000000001
$ws_LTW.Range("L132").Value = 308577.9
$ws_LTW.Range("M132").Value = -71451.20000000001
$ws_LTW.Range("N132").Value = -313637.9

# --- WVR ---
$ws_WVR.Range("H107").Value = 544.2
$ws_WVR.Range("I107").Value = 446.3
$ws_WVR.Range("J107").Value = 740
$ws_WVR.Range("K107").Value = 1338.9
$ws_WVR.Range("L107").Value = 2220
$ws_WVR.Range("M107").Value = 581.0999999999999
$ws_WVR.Range("N107").Value = -6060
$ws_WVR.Range("H113").Value = 333.08334
$ws_WVR.Range("J113").Value = 430.75
$ws_WVR.Range("L113").Value = 1292.25
$ws_WVR.Range("N113").Value = -5632.25
$ws_WVR.Range("H124").Value = 30619.334
$ws_WVR.Range("J124").Value = 30619.334
$ws_WVR.Range("L124").Value = 30619.334
$ws_WVR.Range("N124").Value = -40439.334
$ws_WVR.Range("H132").Value = 2767.0588
$ws_WVR.Range("I132").Value = 3065.4
$ws_WVR.Range("J132").Value = 2340.8572
$ws_WVR.Range("K132").Value = 9196.200000000001
$ws_WVR.Range("L132").Value = 7022.571599999999
$ws_WVR.Range("M132").Value = -6666.200000000001
$ws_WVR.Range("N132").Value = -12082.5716
